# Auto-generated edit script: update crypto price/volume table
# (mirrors the Wed Oct 25 12:50:42 UTC 2023 GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.306.23'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.784.93'
$ws.Range("E3").Value = '  -3.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.66'
$ws.Range("E5").Value = '  -3.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  -4.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '33.58'
$ws.Range("E8").Value = '  +4.58%  '
$ws.Range("E9").Value = '  -3.65%  '
$ws.Range("E10").Value = '  -5.26%  '
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("D12").Value = '2.041.16'
$ws.Range("E12").Value = '  -3.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  +6.57%  '
$ws.Range("D14").Value = '1.799.94'
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.632'
$ws.Range("E15").Value = '  -4.57%  '
$ws.Range("D16").Value = '34.290.26'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.93'
$ws.Range("E18").Value = '  -3.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '255.09'
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").Value = '0.0₃0742'
$ws.Range("E20").Value = '  -3.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.44'
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("E23").Value = '  -6.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("E24").Value = '  -5.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.56'
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("E26").Value = '  -3.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.02'
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("E28").Value = '  -5.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").Value = '  -3.51%  '
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("D35").Value = '1.446.88'
$ws.Range("E35").Value = '  -7.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.624'
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.85'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.98'
$ws.Range("E40").Value = '  -4.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.35'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.890'
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.08'
$ws.Range("E43").Value = '  -4.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0508'
$ws.Range("E44").Value = '  -3.84%  '
$ws.Range("E45").Value = '  -2.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.34'
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.85'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.941.68'
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.51'
$ws.Range("E50").Value = '  -1.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.94'
$ws.Range("E51").Value = '  -4.94%  '
